$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "YATIKA JENA"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "230102105"

$ws.Range("C4").Value = "Very Poor"
$ws.Range("D4").Value = "Very Poor"
$ws.Range("E4").Value = "Very Poor"
$ws.Range("F4").Value = "hojaaa"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "5/17/2025"

$ws.Range("H4").Value = "Unknown"
$ws.Range("I4").Value = "Unknown"
